$d = $word.ActiveDocument

# 1) Update the second "Working 2logLR" line: 19.7 -> 32.7, p= 0.1 -> p= 0.05
$d.Content.Find.Execute(
    "Working 2logLR =  19.7 p= 0.1 ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Working 2logLR =  32.7 p= 0.05 ",
    2)

# 2) Update the scale-factors line (expanded list of values) and the following
#    ";  denominator " text stays the same, only the bracketed numbers change.
$d.Content.Find.Execute(
    "(scale factors:  1.9 1.6 1.3 1 1 0.99 0.91 0.84 0.84 0.82 0.79 0.78 0.74 0.43 );  denominator ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(scale factors:  1.9 1.6 1.5 1.2 1.2 1.1 1 0.98 0.96 0.96 0.93 0.92 0.88 0.85 0.83 0.81 0.77 0.75 0.58 0.27 );  denominator ",
    2)

# 3) Update the final denominator df value: 2107 -> 2101
$d.Content.Find.Execute(
    "= 2107",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "= 2101",
    2)
